# Apply updated crypto price/volume data per the Oct 17 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.962.05'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '2.591.64'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('D5').Formula = "'588.29"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.14%  '
$ws.Range('D6').Formula = "'148.97"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.63%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Formula = "'0.543"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('D9').Value = '2.589.16'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E10').Value = '  -3.54%  '
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('D12').Formula = "'5.15"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('E13').Value = '  -3.33%  '
$ws.Range('D14').Formula = "'27.05"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.19%  '
$ws.Range('D15').Value = '3.062.39'
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('E16').Value = '  -4.93%  '
$ws.Range('D17').Value = '66.959.65'
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D18').Value = '2.591.81'
$ws.Range('E18').Value = '  -1.48%  '
$ws.Range('D19').Formula = "'361.12"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.34%  '
$ws.Range('D20').Formula = "'10.94"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.23%  '
$ws.Range('D21').Formula = "'7.29"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.57%  '
$ws.Range('E22').Value = '  -0.58%  '
$ws.Range('D23').Formula = "'4.81"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.53%  '
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('D25').Formula = "'72.04"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.46%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').Formula = "'9.90"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('D29').Formula = "'1.00"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Formula = "'573.27"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('D31').Value = '0.0₃0974'
$ws.Range('E31').Value = '  -6.42%  '
$ws.Range('D32').Formula = "'1.36"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.95%  '
$ws.Range('D33').Formula = "'7.57"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.38%  '
$ws.Range('E34').Value = '  -3.74%  '
$ws.Range('D35').Formula = "'0.999"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E36').Value = '  -6.15%  '
$ws.Range('E37').Value = '  -3.01%  '
$ws.Range('D38').Formula = "'156.33"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('D39').Formula = "'18.84"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.62%  '
$ws.Range('D40').Formula = "'0.364"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('E41').Value = '  -0.94%  '
$ws.Range('D42').Formula = "'5.15"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.62%  '
$ws.Range('E43').Value = '  +2.05%  '
$ws.Range('D44').Formula = "'2.48"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.95%  '
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').Formula = "'151.53"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.53%  '
$ws.Range('D47').Value = '0.0₆0282'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('D48').Formula = "'3.70"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.55%  '
$ws.Range('B49').Value = 'Optimism'
$ws.Range('C49').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D49').Formula = "'1.67"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.93%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Formula = "'0.0776"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.79%  '
$ws.Range('D51').Formula = "'21.20"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.97%  '
